# Applies the "Added final Avihepadna EVE set" edit:
#  - row 47 (ehbv-avi.34-leptosomus) keeps its place but gets new
#    nearest_upstream_orf / nearest_downstream_orf values (F/G)
#  - a brand-new row 48 (ehbv-avi.35-calypte) is inserted, pushing the
#    former rows 48/49 (callipepla / phylloscopus) down to 49/50
#  - ten brand-new rows 51-60 (ehbv-avi.38 .. ehbv-avi.47) are appended
#  - the sheet view scroll/selection state is updated to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Row 47 (ehbv-avi.34-leptosomus): only F/G (the orf columns) change
# ---------------------------------------------------------------------
$ws.Range("F47").Value = "LRRC2"
$ws.Range("G47").Value = "nk"

# ---------------------------------------------------------------------
# 2. Insert a new row 48 for ehbv-avi.35-calypte, shifting the two
#    rows beneath it (callipepla/phylloscopus) down to 49/50.
#    A fresh, unformatted row is pulled in (copied from well below the
#    table) so that column A keeps the *unfilled* look the new rows use
#    in the workbook (no yellow highlight), then columns B:K pick up
#    the normal "data row" formatting from an existing data row.
# ---------------------------------------------------------------------
$ws.Rows("48").Insert()
$ws.Range("A200:K200").Copy()
$ws.Range("A48:K48").PasteSpecial(-4122)
$ws.Range("B40:K40").Copy()
$ws.Range("B48:K48").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A48").Value = "ehbv-avi.35-calypte"
$ws.Range("B48").Value = "Avihepadnavirus"
$ws.Range("C48").Value = "Endogenous avihepadnavirus 35"
$ws.Range("D48").Value = "ehbv-avi.35-calypte"
$ws.Range("E48").Value = 35
$ws.Range("F48").Value = "TIMM21/FBXO15"
$ws.Range("G48").Value = "CYB5A"
$ws.Range("H48").Value = "Hepadnaviridae"
$ws.Range("I48").Value = "Calypte"
$ws.Range("J48").Value = "N/A"
$ws.Range("K48").Value = "NK"

# ---------------------------------------------------------------------
# 3. Ten brand-new rows 51-60 (ehbv-avi.38 .. ehbv-avi.47), added after
#    the existing last row (50, ehbv-avi.37-phylloscopus). These rows
#    did not exist before, so we just format + fill them directly.
# ---------------------------------------------------------------------
$ws.Range("A200:K209").Copy()
$ws.Range("A51:K60").PasteSpecial(-4122)
$ws.Range("B40:K40").Copy()
$ws.Range("B51:K60").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$newRows = @(
    @{Row=51; A="ehbv-avi.38-passeriformes";     C="Endogenous avihepadnavirus 38"; E=38; F="AAR2";                G="ENSTGUG00000027480"; I="passeriformes"},
    @{Row=52; A="ehbv-avi.39-ara";                C="Endogenous avihepadnavirus 39"; E=39; F="FXN";                 G="FXN";                I="ara"},
    @{Row=53; A="ehbv-avi.40-oxyura";             C="Endogenous avihepadnavirus 40"; E=40; F="ENSACOG00000002782";  G="ENSACOG00000002891"; I="oxyura"},
    @{Row=54; A="ehbv-avi.41-psittaciformes";     C="Endogenous avihepadnavirus 41"; E=41; F="RAD23B";              G="PHAX/MARCHF3";        I="psittaciformes"},
    @{Row=55; A="ehbv-avi.42-passeriformes-con";  C="Endogenous avihepadnavirus 42"; E=42; F="RAD23B";              G="PHAX/MARCHF3";        I="passeriformes"},
    @{Row=56; A="ehbv-avi.43-gallirallus";        C="Endogenous avihepadnavirus 43"; E=43; F="ABRACL/REPS1";       G="TXLNB";               I="gallirallus"},
    @{Row=57; A="ehbv-avi.44-antrostomus";        C="Endogenous avihepadnavirus 44"; E=44; F="ENSMUNG00000008889"; G="ENSMUNG00000008889"; I="antrostomus"},
    @{Row=58; A="ehbv-avi.45-ara";                C="Endogenous avihepadnavirus 45"; E=45; F="KCNV1";               G="ENSTGUG00000027711"; I="ara"},
    @{Row=59; A="ehbv-avi.46-psittaciformes";     C="Endogenous avihepadnavirus 46"; E=46; F="LNPEP";               G="ENSACOG00000006853"; I="psittaciformes"},
    @{Row=60; A="ehbv-avi.47-passer";             C="Endogenous avihepadnavirus 47"; E=47; F="RYR3";                G="FMN1";                I="passer"}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = "Avihepadnavirus"
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.A
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = "Hepadnaviridae"
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = "N/A"
    $ws.Range("K$row").Value = "NK"
}

# ---------------------------------------------------------------------
# 4. Update the sheet view to match (selection)
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("G53").Select()
